$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header labels reordered
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "bedrooms_2"
$ws.Range("E1").Value = "living_rooms_1"
$ws.Range("F1").Value = "living_rooms_2"

# Row 3: B3 0->1, C3 1->0
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0

# Row 4: E4 1->0, F4 0->1
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1

# Row 5: C5 0->1, F5 1->0
$ws.Range("C5").Value = 1
$ws.Range("F5").Value = 0

# Row 6: A6 0->1, B6 1->0
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 0

# Row 7: A7 1->0, E7 0->1
$ws.Range("A7").Value = 0
$ws.Range("E7").Value = 1
